$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New defined name "Fc" -> Sheet1!$B$25 (added before other edits so the
# Fc-dependent formula below resolves instead of erroring out)
$wb.Names.Add("Fc", "=Sheet1!`$B`$25")

# Re-apply the existing G6:G15 / H6:H15 formulas across the whole column so
# the engine folds them into shared formulas (t="shared"), matching how
# Excel stores a fill-down/copy-paste of an identical formula.
$ws.Range("G6:G15").Formula = "=F6*Multi"
$ws.Range("H6:H15").Formula = "=G6*tic_lng"

# New "Antenna calculations" block, rows 24-27
$ws.Range("A24").Value = "Antenna calculations"
$ws.Range("A24").Style = "Accent1"

$ws.Range("A25").Value = "Fcarrier"
$ws.Range("B25").Value = 315
$ws.Range("C25").Value = "MHz"

$ws.Range("A26").Value = "Wave length"
$ws.Range("B26").Formula = "=3*10^8/(Fc*10^6)"
$ws.Range("B26").NumberFormat = "0.00"
$ws.Range("C26").Value = "m"

$ws.Range("A27").Value = "WL/4"
$ws.Range("B27").Formula = "=B26/4"
$ws.Range("B27").NumberFormat = "0.00"
$ws.Range("C27").Value = "m"

# Move the active selection like the source workbook (was D21, now D29)
$ws.Range("D29").Select() | Out-Null
